# Moved "Interesting Counter-Arguments" sheet's content out (deleted the sheet)
# and reworked the Venue breakdown table on the Analysis sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Remove the "Interesting Counter-Arguments" worksheet entirely.
# ---------------------------------------------------------------------------
$wsCounter = $wb.Worksheets.Item("Interesting Counter-Arguments")
$wsCounter.Delete()

# ---------------------------------------------------------------------------
# 2. Rework the "Venue" breakdown table on the Analysis sheet.
#    Old layout (rows 37-46): JSSC, ISSCC, VLSI, CICC, ESSCIRC, ASSCC, ISCAS,
#    TCAS-I, TCAS-II, Other
#    New layout (rows 38-49): ISSCC, VLSI, CICC, ESSERC(+ESSCIRC), ASSCC,
#    RFIC, ISCAS, JSSC, TCAS-I, TCAS-II, TVLSI, Other
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Analysis")

# Make room: push the existing table down by one row (37 -> 38 .. 46 -> 47)
$ws2.Rows.Item(37).Insert()
# ... and that leaves a blank spacer row where the table header used to
# transition directly into data - clear it so it disappears like the other
# blank spacer rows in this sheet.
$ws2.Rows.Item(37).ClearContents()
$ws2.Rows.Item(37).ClearFormats()

# Add two more rows at the bottom of the table (for the two new venues).
$ws2.Rows.Item(48).Insert()
$ws2.Rows.Item(49).Insert()

# Now rows 38-49 are available (38-47 hold the old table data, 48-49 blank).
# Overwrite the whole block with the final values/formulas.
$ws2.Range("A38").Value = "ISSCC"
$ws2.Range("B38").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A38)"

$ws2.Range("A39").Value = "VLSI"
$ws2.Range("B39").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A39)"

$ws2.Range("A40").Value = "CICC"
$ws2.Range("B40").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A40)"

$ws2.Range("A41").Value = "ESSERC"
$ws2.Range("B41").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,""ESSCIRC"")+COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,""ESSERC"")"

$ws2.Range("A42").Value = "ASSCC"
$ws2.Range("B42").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A42)"

$ws2.Range("A43").Value = "RFIC"
$ws2.Range("B43").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A43)"

$ws2.Range("A44").Value = "ISCAS"
$ws2.Range("B44").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A44)"

$ws2.Range("A45").Value = "JSSC"
$ws2.Range("B45").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A45)"

$ws2.Range("A46").Value = "TCAS-I"
$ws2.Range("B46").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A46)"

$ws2.Range("A47").Value = "TCAS-II"
$ws2.Range("B47").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A47)"

$ws2.Range("A48").Value = "TVLSI"
$ws2.Range("B48").Formula = "=COUNTIF('Ringamp Publication List'!`$B`$2:`$B`$997,Analysis!A48)"

$ws2.Range("A49").Value = "Other"
$ws2.Range("B49").Formula = "=COUNTA('Ringamp Publication List'!`$A`$2:`$A`$997)-SUM(B38:B48)"

# ---------------------------------------------------------------------------
# 3. Point the "Publications by Venue" chart at the new range and let it
#    recompute its cache from the refreshed cells.
# ---------------------------------------------------------------------------
$chartObjects = $ws2.ChartObjects()
$venueChart = $chartObjects.Item(3).Chart
$venueSeries = $venueChart.SeriesCollection(1)
$venueSeries.Formula = "=SERIES(,Analysis!`$A`$38:`$A`$49,Analysis!`$B`$38:`$B`$49,1)"

# ---------------------------------------------------------------------------
# 4. Window/view bookkeeping: Analysis becomes the active/selected sheet.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("M49").Select()
